$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: finished 84 pages of Arooj-e-Iqbal (Azka), logged via the
# existing "Book Title" / "Time Duration" shared strings.
$ws.Range("C5").Value = "Arooj-e-Iqbal"
$ws.Range("D5").Value = 120
$ws.Range("E5").Value = 201
$ws.Range("F5").Value = 84

$ws.Range("G5").Value = 0.32291666666666669
$ws.Range("G5").NumberFormat = "h:mm"

$ws.Range("H5").Value = 0.88194444444444453
$ws.Range("H5").NumberFormat = "h:mm"

$ws.Range("I5").Value = "1.5 hours"
$ws.Range("J5").Value = 84

# Match the row height used by the rest of the sheet.
$ws.Rows(5).RowHeight = 15.75

# Move the selection the way it ended up after the edit.
[void]$ws.Range("J6").Select()
